$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.040535926818848
$ws.Range("B1").Value = 3.473752498626709
$ws.Range("C1").Value = 3.473796129226685
$ws.Range("D1").Value = 2.028436183929443
$ws.Range("E1").Value = 1.167825698852539
